# Add files via upload
# - Append a new row (1309, 1309, "ansible-roles", "ansible-roles", "ansible-roles")
#   to the "groups" worksheet.
# - Make "groups" the active sheet/tab, with selection E17.
# - Leave "projects" sheet selection at D2 and no longer the active tab.

$wb = $excel.ActiveWorkbook

$wsGroups = $wb.Worksheets.Item("groups")
$wsProjects = $wb.Worksheets.Item("projects")

# Append new data row to "groups" sheet
$wsGroups.Range("A14").Value = 1309
$wsGroups.Range("B14").Value = 1309
$wsGroups.Range("C14").Value = "ansible-roles"
$wsGroups.Range("D14").Value = "ansible-roles"
$wsGroups.Range("E14").Value = "ansible-roles"

# Update selection on the "projects" sheet before switching away from it
$wsProjects.Range("D2").Select()

# Activate "groups" sheet and set its selection
$wsGroups.Activate()
$wsGroups.Range("E17").Select()
